# Edited last slide of presentation. Added pdf.
#
# The title placeholder on the final "Thank you" slide is moved/resized
# (made wider and nudged down) and its text is expanded from
# "Благодаря ви!" to "Благодаря за вниманието".

$p = $ppt.ActivePresentation

# Last slide of the deck (the "Благодаря ..." / thank-you slide).
$lastSlide = $p.Slides.Item($p.Slides.Count)

# "Title 1" placeholder is the first shape on that slide.
$title = $lastSlide.Shapes.Item(1)

# New position/size for the title placeholder (values chosen so the
# EMU-rounded result lands exactly on the target offsets/extents).
$title.Left   = 24.230629921259844
$title.Top    = 207.92308044433594
$title.Width  = 565.6155118110236
$title.Height = 106.8

# Update the title text.
$title.TextFrame.TextRange.Text = "Благодаря за вниманието"
